$p = $ppt.ActivePresentation

# --- Slide 4: "Other questions" -------------------------------------------
# Insert four new bullet questions ahead of the existing one.
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange

$lines4 = @(
    "Does using an AR1 spatiotemporal field (without factor levels for years) constrain the model too much and result in hyperstability?",
    "Are there survey designs that result in the model-based index being more or less precise than the design based index?",
    "If the catchability of a survey changed along the time series, say the gear was changed and there was one year of calibration overlap, could the model estimate the catchability (q) offset and provide unbiased estimates of the population available to the contemporary survey as if those gear were used the whole time?",
    "Can we obtain an index at age using a geostatistical model?",
    "Does the model sufficiently account for spatial correlation and/or is it sufficiently free of assumption so as not to be affected by the stratified sampling design of the survey?"
)
$tr4.Text = [string]::Join("`r", $lines4)

# --- Slide 5: "Approach taken" ---------------------------------------------
# Replace the single bullet with a numbered list describing the approach.
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(2)
$tr5 = $shp5.TextFrame.TextRange

$lines5 = @(
    "Simulate a population and a survey and calculate design-based indices using SimSurvey",
    "Fit a geostistical model to the simulated survey data using sdmTMB to obtain model-based indices",
    "Iterate the population simulation and data analysis",
    "Visually assess the bias and precision of the estimates",
    "Modify the simulation settings (e.g., impose partial survey coverage) and repeat setps 1-4"
)
$tr5.Text = [string]::Join("`r", $lines5)

for ($i = 1; $i -le $lines5.Count; $i++) {
    $para = $tr5.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Type = 2
}
